$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B126").Value = 65258
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("B127").Value = 64196
$ws.Range("F127").Value = 1
$ws.Range("G127").Value = 32143.58
$ws.Range("B161").Value = 64350
$ws.Range("E161").Value = 70.63
$ws.Range("F161").Value = 2
$ws.Range("G161").Value = 132.88
$ws.Range("B162").Value = 57756
$ws.Range("E162").Value = 79.37
$ws.Range("F162").Value = -100
$ws.Range("G162").Value = -6644
$ws.Range("F177").Value = 279
$ws.Range("G177").Value = 12993.03
$ws.Range("B193").Value = 68046.69
$ws.Range("F203").Value = 6
$ws.Range("G203").Value = 50.7
$ws.Range("B204").Value = 50.7
$ws.Range("F206").Value = 75
$ws.Range("G206").Value = 4860
$ws.Range("B208").Value = 4906.49
$ws.Range("F222").Value = 960
$ws.Range("G222").Value = 17760
$ws.Range("B229").Value = 29994.07
$ws.Range("F263").Value = 17
$ws.Range("G263").Value = 1762.9
$ws.Range("F267").Value = 139
$ws.Range("G267").Value = 5904.72
$ws.Range("B292").Value = 66196
$ws.Range("C292").Value = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F292").Value = 6
$ws.Range("G292").Value = 526.2
$ws.Range("B293").Value = 64985
$ws.Range("C293").Value = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("F293").Value = 12
$ws.Range("G293").Value = 1052.4
$ws.Range("B295").Value = 125806.45
$ws.Range("B308").Value = 55356
$ws.Range("E308").Value = 54.04
$ws.Range("F308").Value = -158
$ws.Range("G308").Value = -7527.12
$ws.Range("B309").Value = 63510
$ws.Range("E309").Value = 50.66
$ws.Range("F309").Value = 78
$ws.Range("G309").Value = 3715.92
$ws.Range("B317").Value = 63560
$ws.Range("E317").Value = 134.87
$ws.Range("F317").Value = 1
$ws.Range("G317").Value = 126.86
$ws.Range("B318").Value = 60325
$ws.Range("E318").Value = 151.57
$ws.Range("F318").Value = -102
$ws.Range("G318").Value = -12939.72
$ws.Range("F325").Value = 41
$ws.Range("G325").Value = 6197.97
$ws.Range("B328").Value = -2344.99
$ws.Range("F361").Value = 248
$ws.Range("G361").Value = 34866.32
$ws.Range("B363").Value = 79376.33
$ws.Range("F387").Value = 439
$ws.Range("G387").Value = 42407.4
$ws.Range("B389").Value = 59348.86
$ws.Range("F402").Value = 51
$ws.Range("G402").Value = 1749.81
$ws.Range("F403").Value = 74
$ws.Range("G403").Value = 2999.96
$ws.Range("B417").Value = 174436.32
$ws.Range("F452").Value = 56
$ws.Range("G452").Value = 15188.32
$ws.Range("B458").Value = 100763.31
$ws.Range("F465").Value = 27
$ws.Range("G465").Value = 3140.64
$ws.Range("F466").Value = 15
$ws.Range("G466").Value = 618
$ws.Range("F469").Value = 2
$ws.Range("G469").Value = 111.78
$ws.Range("B476").Value = 49028.15
$ws.Range("B479").Value = 53319
$ws.Range("E479").Value = 310.64
$ws.Range("F479").Value = -6
$ws.Range("G479").Value = -1643.52
$ws.Range("B480").Value = 64810
$ws.Range("E480").Value = 291.22
$ws.Range("F480").Value = 0
$ws.Range("G480").Value = 0
$ws.Range("B496").Value = 60025
$ws.Range("E496").Value = 37.22
$ws.Range("F496").Value = -98
$ws.Range("G496").Value = -3217.34
$ws.Range("B497").Value = 64833
$ws.Range("E497").Value = 34.9
$ws.Range("F497").Value = 88
$ws.Range("G497").Value = 2889.04
$ws.Range("B506").Value = 64830
$ws.Range("E506").Value = 34.9
$ws.Range("F506").Value = 85
$ws.Range("G506").Value = 2790.55
$ws.Range("B507").Value = 60022
$ws.Range("E507").Value = 37.22
$ws.Range("F507").Value = -113
$ws.Range("G507").Value = -3709.79
$ws.Range("F519").Value = 422
$ws.Range("G519").Value = 23159.36
$ws.Range("B525").Value = 130450.51
$ws.Range("F530").Value = 28
$ws.Range("G530").Value = 1209.04
$ws.Range("F532").Value = 11
$ws.Range("G532").Value = 474.98
$ws.Range("F533").Value = 0
$ws.Range("G533").Value = 0
$ws.Range("B535").Value = 25706.38
$ws.Range("F605").Value = 191
$ws.Range("G605").Value = 25422.1
$ws.Range("F606").Value = 3
$ws.Range("G606").Value = 405.03
$ws.Range("B607").Value = 25827.13
$ws.Range("F609").Value = 23
$ws.Range("G609").Value = 2502.63
$ws.Range("F612").Value = 236
$ws.Range("G612").Value = 35496.76
$ws.Range("F615").Value = 103
$ws.Range("G615").Value = 15931.01
$ws.Range("F625").Value = 330
$ws.Range("G625").Value = 12153.9
$ws.Range("B628").Value = 215925.92
$ws.Range("F662").Value = 48
$ws.Range("G662").Value = 3854.88
$ws.Range("B668").Value = 13039.27
$ws.Range("F674").Value = 906
$ws.Range("G674").Value = 147777.66
$ws.Range("B680").Value = 148790.21
$ws.Range("B718").Value = 2835067.44
$ws.Range("B719").Value = 2835067.44
